$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 3 and 4
$wsOverview.Range("G3").Value = "2016-08-18 20:16:04"
$wsOverview.Range("G4").Value = "2016-08-18 20:16:04"

# zh-cn sheet: Priority column (E) rows 3 & 4: ht -> mt
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime column (H) rows 3 & 4
$wsZhCn.Range("H3").Value = "2016-08-18 20:15:55"
$wsZhCn.Range("H4").Value = "2016-08-18 20:15:55"

# zh-cn sheet: Correspond Handback DateTime column (K) rows 3 & 4
$wsZhCn.Range("K3").Value = "2016-08-18 20:16:26"
$wsZhCn.Range("K4").Value = "2016-08-18 20:16:26"

# de-de sheet: Correspond Handoff Datetime column (H) rows 3 & 4 (shares string with Overview G3/G4)
$wsDeDe.Range("H3").Value = "2016-08-18 20:16:04"
$wsDeDe.Range("H4").Value = "2016-08-18 20:16:04"

# de-de sheet: Priority column (E) rows 3 & 4: ht -> mt (shares string with zh-cn E3/E4)
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handback DateTime column (K) rows 3 & 4
$wsDeDe.Range("K3").Value = "2016-08-18 20:16:35"
$wsDeDe.Range("K4").Value = "2016-08-18 20:16:35"
